# Removed frequencies from the test automation configuration xlsx
#
# The original sheet had a block of simulation "frequency" rows
# (GnssPoseSimulink, PointsRawFloat32, ImageRaw, ClockFrequency,
# SimulinkState, CurrentVelocity, PoseOtherCar, CurrentPose) between the
# GoalOrientW row and the percent_reflecting_sfc/R rows. This change
# deletes that whole block of rows, shifting percent_reflecting_sfc/R up,
# and sets percent_reflecting_sfc's value to 0.9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 19-26 hold the 8 now-unwanted frequency entries. Deleting them
# shifts the following rows (percent_reflecting_sfc, R) up to rows 19-20.
$ws.Range("A19:B26").EntireRow.Delete() | Out-Null

# percent_reflecting_sfc (now row 19) gets a new value.
$ws.Range("B19").Value = 0.9

# Cosmetic: keep view state close to what Excel would leave behind after
# selecting/deleting that block of rows.
try {
    $excel.ActiveWindow.ScrollRow = 7
    $excel.ActiveWindow.TopLeftCell = $ws.Range("A7")
} catch {
}
$ws.Range("A19:XFD26").Select() | Out-Null
